# Auto-generated edit script: refreshes market-price-derived columns
# (currentAveragePrice*, LevePrice*, LeveProfit*) on each job sheet,
# matching the scheduled market-data refresh described in the commit.
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 462.25
$ws.Range("J19").Value = 483
$ws.Range("L19").Value = 483
$ws.Range("N19").Value = -833
$ws.Range("H28").Value = 877.1905
$ws.Range("I28").Value = 965.94116
$ws.Range("K28").Value = 965.94116
$ws.Range("M28").Value = -480.94116
$ws.Range("H70").Value = 9246.27
$ws.Range("I70").Value = 10812.143
$ws.Range("J70").Value = 7419.4165
$ws.Range("K70").Value = 32436.429
$ws.Range("L70").Value = 22258.2495
$ws.Range("M70").Value = -32166.429
$ws.Range("N70").Value = -22798.2495
$ws.Range("H73").Value = 9246.27
$ws.Range("I73").Value = 10812.143
$ws.Range("J73").Value = 7419.4165
$ws.Range("K73").Value = 32436.429
$ws.Range("L73").Value = 22258.2495
$ws.Range("M73").Value = -31500.429
$ws.Range("N73").Value = -24130.2495
$ws.Range("H98").Value = 1560
$ws.Range("I98").Value = 1331
$ws.Range("K98").Value = 1331
$ws.Range("M98").Value = 167
$ws.Range("H116").Value = 9004.944
$ws.Range("I116").Value = 13694.111
$ws.Range("K116").Value = 13694.111
$ws.Range("M116").Value = -10252.111
$ws.Range("H122").Value = 1560
$ws.Range("I122").Value = 1331
$ws.Range("K122").Value = 3993
$ws.Range("M122").Value = -1543
$ws.Range("H139").Value = 154780
$ws.Range("J139").Value = 154780
$ws.Range("L139").Value = 154780
$ws.Range("N139").Value = -165060

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1540
$ws.Range("I2").Value = 1425.7693
$ws.Range("K2").Value = 1425.7693
$ws.Range("M2").Value = -1312.7693
$ws.Range("H32").Value = 4980.7334
$ws.Range("I32").Value = 4106.327
$ws.Range("K32").Value = 4106.327
$ws.Range("M32").Value = -3819.327
$ws.Range("H61").Value = 1774.5807
$ws.Range("I61").Value = 1774.5807
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 1774.5807
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -1562.5807
$ws.Range("N61").ClearContents()
$ws.Range("H74").Value = 1036.825
$ws.Range("I74").Value = 919.8158
$ws.Range("K74").Value = 919.8158
$ws.Range("M74").Value = -45.81579999999997
$ws.Range("H77").Value = 1036.825
$ws.Range("I77").Value = 919.8158
$ws.Range("K77").Value = 4599.079
$ws.Range("M77").Value = -231.0789999999997
$ws.Range("H97").Value = 9950.166999999999
$ws.Range("I97").Value = 3622.5
$ws.Range("J97").Value = 22605.5
$ws.Range("K97").Value = 3622.5
$ws.Range("L97").Value = 22605.5
$ws.Range("M97").Value = -3126.5
$ws.Range("N97").Value = -23597.5
$ws.Range("H110").Value = 3184.6875
$ws.Range("I110").Value = 1593.1
$ws.Range("K110").Value = 1593.1
$ws.Range("M110").Value = 451.9000000000001
$ws.Range("H116").Value = 1540
$ws.Range("I116").Value = 1425.7693
$ws.Range("K116").Value = 1425.7693
$ws.Range("M116").Value = 868.2307000000001
$ws.Range("H122").Value = 1936.3489
$ws.Range("I122").Value = 1564.7037
$ws.Range("J122").Value = 2563.5
$ws.Range("K122").Value = 4694.1111
$ws.Range("L122").Value = 7690.5
$ws.Range("M122").Value = -2244.1111
$ws.Range("N122").Value = -12590.5
$ws.Range("H136").Value = 1774.5807
$ws.Range("I136").Value = 1774.5807
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 5323.742099999999
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -2773.742099999999
$ws.Range("N136").ClearContents()

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1540
$ws.Range("I3").Value = 1425.7693
$ws.Range("K3").Value = 1425.7693
$ws.Range("M3").Value = -1311.7693
$ws.Range("H53").Value = 99999
$ws.Range("J53").Value = 99999
$ws.Range("L53").Value = 99999
$ws.Range("N53").Value = -101147
$ws.Range("H132").Value = 134567
$ws.Range("J132").Value = 134567
$ws.Range("L132").Value = 134567
$ws.Range("N132").Value = -144687
$ws.Range("H134").Value = 1872.0625
$ws.Range("I134").Value = 1496.3846
$ws.Range("K134").Value = 4489.1538
$ws.Range("M134").Value = -1954.1538

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 3484.182
$ws.Range("I86").Value = 2354.1667
$ws.Range("K86").Value = 2354.1667
$ws.Range("M86").Value = -1231.1667
$ws.Range("H89").Value = 3484.182
$ws.Range("I89").Value = 2354.1667
$ws.Range("K89").Value = 11770.8335
$ws.Range("M89").Value = -6154.833500000001
$ws.Range("H99").Value = 7869.6523
$ws.Range("I99").Value = 11546
$ws.Range("J99").Value = 3090.4
$ws.Range("K99").Value = 11546
$ws.Range("L99").Value = 3090.4
$ws.Range("M99").Value = -10048
$ws.Range("N99").Value = -6086.4
$ws.Range("H126").Value = 7869.6523
$ws.Range("I126").Value = 11546
$ws.Range("J126").Value = 3090.4
$ws.Range("K126").Value = 34638
$ws.Range("L126").Value = 9271.200000000001
$ws.Range("M126").Value = -32168
$ws.Range("N126").Value = -14211.2
$ws.Range("H132").Value = 4498.048
$ws.Range("I132").Value = 5250.1333
$ws.Range("J132").Value = 2617.8333
$ws.Range("K132").Value = 15750.3999
$ws.Range("L132").Value = 7853.499899999999
$ws.Range("M132").Value = -13220.3999
$ws.Range("N132").Value = -12913.4999
$ws.Range("H141").Value = 299746.5
$ws.Range("J141").Value = 299746.5
$ws.Range("L141").Value = 299746.5
$ws.Range("N141").Value = -310106.5

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H80").Value = 3162.8333
$ws.Range("I80").Value = 1803.5714
$ws.Range("J80").Value = 3722.5293
$ws.Range("K80").Value = 5410.7142
$ws.Range("L80").Value = 11167.5879
$ws.Range("M80").Value = -4474.7142
$ws.Range("N80").Value = -13039.5879
$ws.Range("H83").Value = 3162.8333
$ws.Range("I83").Value = 1803.5714
$ws.Range("J83").Value = 3722.5293
$ws.Range("K83").Value = 16232.1426
$ws.Range("L83").Value = 33502.7637
$ws.Range("M83").Value = -11552.1426
$ws.Range("N83").Value = -42862.7637
$ws.Range("H107").Value = 590.7646999999999
$ws.Range("I107").Value = 642.5
$ws.Range("J107").Value = 583.86664
$ws.Range("K107").Value = 1927.5
$ws.Range("L107").Value = 1751.59992
$ws.Range("M107").Value = -7.5
$ws.Range("N107").Value = -5591.59992
$ws.Range("H113").Value = 2279.0588
$ws.Range("J113").Value = 2279.0588
$ws.Range("L113").Value = 6837.176399999999
$ws.Range("N113").Value = -11177.1764
$ws.Range("H124").Value = 0
$ws.Range("I124").Value = 0
$ws.Range("K124").Value = 0
$ws.Range("M124").ClearContents()

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H57").Value = 16016.333
$ws.Range("I57").Value = 9025
$ws.Range("K57").Value = 9025
$ws.Range("M57").Value = -8205
$ws.Range("H97").Value = 27794.45
$ws.Range("I97").Value = 41381.21
$ws.Range("J97").Value = 1979.6
$ws.Range("K97").Value = 41381.21
$ws.Range("L97").Value = 1979.6
$ws.Range("M97").Value = -40885.21
$ws.Range("N97").Value = -2971.6
$ws.Range("H102").Value = 1417.0938
$ws.Range("I102").Value = 1461.6666
$ws.Range("K102").Value = 1461.6666
$ws.Range("M102").Value = 160.3334
$ws.Range("H113").Value = 3688.0588
$ws.Range("I113").Value = 2819.5833
$ws.Range("K113").Value = 2819.5833
$ws.Range("M113").Value = -649.5832999999998
$ws.Range("H122").Value = 2518.92
$ws.Range("I122").Value = 2881.3125
$ws.Range("K122").Value = 8643.9375
$ws.Range("M122").Value = -6193.9375
$ws.Range("H126").Value = 2601.6
$ws.Range("I126").Value = 2740.6667
$ws.Range("J126").Value = 2393
$ws.Range("K126").Value = 8222.000100000001
$ws.Range("L126").Value = 7179
$ws.Range("M126").Value = -5752.000100000001
$ws.Range("N126").Value = -12119
$ws.Range("H132").Value = 4737.353
$ws.Range("I132").Value = 5317.25
$ws.Range("K132").Value = 15951.75
$ws.Range("M132").Value = -13421.75

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 121077.07
$ws.Range("I61").Value = 102210.8
$ws.Range("K61").Value = 102210.8
$ws.Range("M61").Value = -102008.8
$ws.Range("H68").Value = 1400
$ws.Range("I68").Value = 1400
$ws.Range("K68").Value = 1400
$ws.Range("M68").Value = -651
$ws.Range("H71").Value = 1400
$ws.Range("I71").Value = 1400
$ws.Range("K71").Value = 7000
$ws.Range("M71").Value = -3256
$ws.Range("H100").Value = 21230.55
$ws.Range("I100").Value = 4951.75
$ws.Range("J100").Value = 32083.084
$ws.Range("K100").Value = 4951.75
$ws.Range("L100").Value = 32083.084
$ws.Range("M100").Value = -4410.75
$ws.Range("N100").Value = -33165.084
$ws.Range("H113").Value = 121077.07
$ws.Range("I113").Value = 102210.8
$ws.Range("K113").Value = 102210.8
$ws.Range("M113").Value = -100040.8
$ws.Range("H132").Value = 3182.7346
$ws.Range("I132").Value = 2785.425
$ws.Range("J132").Value = 4948.5557
$ws.Range("K132").Value = 8356.275000000001
$ws.Range("L132").Value = 14845.6671
$ws.Range("M132").Value = -5826.275000000001
$ws.Range("N132").Value = -19905.6671

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H3").Value = 15193886
$ws.Range("I3").Value = 27886770
$ws.Range("J3").Value = 2501002
$ws.Range("K3").Value = 27886770
$ws.Range("L3").Value = 2501002
$ws.Range("M3").Value = -27886656
$ws.Range("N3").Value = -2501230
$ws.Range("H122").Value = 1681.7084
$ws.Range("I122").Value = 1555.4706
$ws.Range("K122").Value = 4666.4118
$ws.Range("M122").Value = -2216.4118
$ws.Range("H132").Value = 2886.75
$ws.Range("I132").Value = 2304.6191
$ws.Range("K132").Value = 6913.8573
$ws.Range("M132").Value = -4383.8573
